$d = $word.ActiveDocument

# Locate the run containing "sl." (end of the sentence "... HTML i sl.")
$rng = $d.Content
$found = $rng.Find.Execute("sl.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'sl.'"
}

# $rng now spans exactly "sl." (Start..End). We only want to touch the
# trailing period, so that the preceding text ("i ") and the "sl" text
# stay in their existing, untouched runs (matches the target diff, which
# leaves the run boundaries before "sl" alone).
$periodRange = $d.Range($rng.End - 1, $rng.End)

# A plain Text/Find-replace on this paragraph re-flows (merges) every
# sibling run that shares identical formatting, which would incorrectly
# fold the preceding "i " run into this edit. Toggling a character
# formatting property first forces the engine to split the "." off into
# its own run without touching neighboring runs.
$periodRange.Font.Bold = $true

# Re-acquire the (now isolated) run holding just "." and replace its text.
$periodRange2 = $d.Range($rng.End - 1, $rng.End)
$periodRange2.Text = "ično."

# Restore normal (non-bold) formatting on the newly inserted text so the
# final run formatting matches the surrounding text exactly.
$newRunRange = $d.Range($rng.End - 1, $rng.End - 1 + 5)
$newRunRange.Font.Bold = $false

Write-Output "Done: 'sl.' -> 'slično.'"
